$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above current row 2 (shifts existing data rows 2-18 down to 6-22)
$ws.Rows("2:5").Insert()


# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2023-10-31"
$ws.Range("B2").Value = "큐로셀"
$ws.Range("C2").Value = "미래, 삼성"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2023-11-03"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2023-11-09"
$ws.Range("F2").Value = 32000000
$ws.Range("G2").Value = 1600000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 29800
$ws.Range("J2").Value = 33500
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "169.95 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2023-10-31"
$ws.Range("B3").Value = "메가터치"
$ws.Range("C3").Value = "NH"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2023-11-03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2023-11-10"
$ws.Range("F3").Value = 24960000
$ws.Range("G3").Value = 5200000
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = 3500
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = "-"
$ws.Range("L3").Value = 4800
$ws.Range("M3").Value = "-"
$ws.Range("N3").Value = "-"
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "-"
$ws.Range("Q3").Value = "-"
$ws.Range("R3").Value = "630.752 : 1"
$ws.Range("S3").Value = "-"
$ws.Range("T3").Value = "-"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "2023-10-31"
$ws.Range("B4").Value = "컨텍"
$ws.Range("C4").Value = "대신"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2023-11-03"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2023-11-09"
$ws.Range("F4").Value = 46350000
$ws.Range("G4").Value = 2060000
$ws.Range("H4").Value = "-"
$ws.Range("I4").Value = 20300
$ws.Range("J4").Value = 22500
$ws.Range("K4").Value = "-"
$ws.Range("L4").Value = 22500
$ws.Range("M4").Value = "-"
$ws.Range("N4").Value = "-"
$ws.Range("O4").Value = 2.912621359223301
$ws.Range("P4").Value = "-"
$ws.Range("Q4").Value = "-"
$ws.Range("R4").Value = "9.09 : 1"
$ws.Range("S4").Value = "-"
$ws.Range("T4").Value = "-"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2023-10-31"
$ws.Range("B5").Value = "비아이매트릭스"
$ws.Range("C5").Value = "IBK"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "2023-11-03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2023-11-09"
$ws.Range("F5").Value = 15600000
$ws.Range("G5").Value = 1200000
$ws.Range("H5").Value = "-"
$ws.Range("I5").Value = 9100
$ws.Range("J5").Value = 11000
$ws.Range("K5").Value = "-"
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = "-"
$ws.Range("N5").Value = "-"
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = "-"
$ws.Range("Q5").Value = "-"
$ws.Range("R5").Value = "805.12 : 1"
$ws.Range("S5").Value = "-"
$ws.Range("T5").Value = "-"

# Clean up any inherited/number formatting on the newly inserted rows so the
# exported cells carry no explicit style (matches the rest of the data rows).
$ws.Range("A2:T5").ClearFormats()
